$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('CS2025_00112','Alice',' alice.johnson@example.com',' Technical',' Unable to log in to my account after recent update.'),
    @('CS2025_00113','Amanda',' amanda.king@example.com',' Technical',' I am getting frequent error messages on the desktop application.'),
    @('CS2025_00114','Benjamin',' benjamin.lewis@example.com',' Billing',' I was charged for a subscription I canceled last month.'),
    @('CS2025_00115','Bob',' bob.smith@example.com',' Billing',' Incorrect charge appeared on my credit card statement.'),
    @('CS2025_00116','Catherine',' catherine.martin@example.com',' General',' What security features are included in your software?'),
    @('CS2025_00117','Clara',' clara.martinez@example.com',' General',' Do you offer discounts for bulk purchases?'),
    @('CS2025_00118','Daniel',' daniel.turner@example.com',' Technical',' The system doesn''t recognize my external hard drive.'),
    @('CS2025_00119','David',' david.lee@example.com',' Technical',' Website shows a 404 error when trying to access the support page.'),
    @('CS2025_00120','Emma',' emma.brown@example.com',' Billing',' I need a copy of my last three invoices for reimbursement purposes.'),
    @('CS2025_00121','Frank',' frank.harris@example.com',' General',' What are your customer service operating hours?'),
    @('CS2025_00122','Grace',' grace.wilson@example.com',' Technical',' My mobile app keeps crashing whenever I try to open the settings page.'),
    @('CS2025_00123','Harry',' harry.adams@example.com',' Billing',' I accidentally made a duplicate payment. Can I get a refund?'),
    @('CS2025_00124','Isabella',' isabella.white@example.com',' General',' Can I change my subscription plan mid-cycle without penalties?'),
    @('CS2025_00125','Jack',' jack.davis@example.com',' Technical',' I am experiencing connectivity issues with your cloud service.'),
    @('CS2025_00126','Karen',' karen.mitchell@example.com',' Billing',' My subscription renewal failed, and my account is deactivated.'),
    @('CS2025_00127','Liam',' liam.robinson@example.com',' General',' The software update failed to install on my device.'),
    @('CS2025_00128','Mia',' mia.scott@example.com',' Technical',' Do you have a physical store near my location?'),
    @('CS2025_00129','Noah',' noah.walker@example.com',' Billing',' I need help configuring my email on a new device.'),
    @('CS2025_00130','Olivia',' olivia.young@example.com',' General',' Can I split my annual payment into monthly installments?'),
    @('CS2025_00131','Paul',' paul.baker@example.com',' Technical',' I am unable to reset my password due to an error.'),
    @('CS2025_00132','Quinn',' quinn.moore@example.com',' Billing',' Can you explain the differences between your plans?'),
    @('CS2025_00133','Rachel',' rachel.thomas@example.com',' General',' Why was I charged twice for the same transaction?'),
    @('CS2025_00134','Sophia',' sophia.hernandez@example.com',' Technical',' My printer isn''t connecting to the Wi-Fi network.'),
    @('CS2025_00135','Tom',' tom.wilson@example.com',' Billing',' Is there a user manual for your latest product?'),
    @('CS2025_00136','Uma',' uma.carter@example.com',' General',' The application freezes when I try to upload large files.'),
    @('CS2025_00137','Victor',' victor.evans@example.com',' Technical',' I need assistance updating my payment method for my account.'),
    @('CS2025_00138','Wendy',' wendy.green@example.com',' Billing',' Do you provide training for your enterprise software?'),
    @('CS2025_00139','Xavier',' xavier.hill@example.com',' General',' I can''t access my account because of a two-factor authentication issue.'),
    @('CS2025_00140','Yolanda',' yolanda.hughes@example.com',' Technical',' How can I get a detailed breakdown of my charges for the last six months?'),
    @('CS2025_00141','Zachary',' zachary.perez@example.com',' Billing',' Does your service offer multi - language support?'),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
